$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 9: mark the "Done" action status as "Done.Theory." for the replication
# topic (F9). New shared string #24 created first (below, for E10) so this
# becomes shared string #25.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 10: fill in the previously-empty Actual Completion Date, Action Status
# and Actions cells for "Design and implement a replication strategy..."
# ---------------------------------------------------------------------------
$ws.Range("E10").Value2 = "1.Data replication across globe.`n2.Cost of data replication.(no of RU's/region* number of regions)`n3.Automatic failover`n4.Manual failover.`n5.Configure Consistancy model.`n6.Configure multi region writes."

$ws.Range("F9").Value2 = "Done.Theory."
$ws.Range("F10").Value2 = "Done.Theory."

$ws.Range("D10").Value2 = 45446
$ws.Range("D10").NumberFormat = $ws.Range("C10").NumberFormat

$ws.Rows.Item(10).RowHeight = 100.8

# ---------------------------------------------------------------------------
# Row 11: new topic row - "Optimize query and operation performance in Azure
# Cosmos DB for NoSQL"
# ---------------------------------------------------------------------------
$ws.Range("A11").Value2 = "Optimize query and operation performance in Azure Cosmos DB for NoSQL"
$ws.Range("B11").Value2 = 2
$ws.Range("C11").Value2 = 45447
$ws.Range("C11").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("E11").Value2 = "1.Custom Index Policy`n2.Read heavy Application.`n3.Write Heavy Application.`n4.Index Metrics (suggest if need any index for query)`n5.Measure Query cost.`n6.Integrated Cache (client’s consistency level must be set to session or eventual)"
$ws.Range("F11").Value2 = "Done"

$ws.Rows.Item(11).RowHeight = 115.2

# ---------------------------------------------------------------------------
# Sheet view: scroll so row 9 is at the top and select F12
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("F12").Select()
